$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# Fill the previously-empty "Definition" column (D) with the same text
# as the "Display" column (C) for each concept row.
$ws.Range("D2").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("C3").Value()
$ws.Range("D4").Value = $ws.Range("C4").Value()
